$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicate "invalidUser@vrbank2" data in row 4: clear A4/B4 contents.
# (A4 keeps its hyperlink-style formatting, matching the target sheet — only the
# value and the B4 cell are removed; the shared string becomes unreferenced and
# is dropped automatically on save.)
$ws.Range("B4").ClearContents()
$ws.Range("A4").ClearContents()

# Remove the individual hyperlink that was anchored on A4 alone (rId4), while
# keeping the A2, A3:A4 (range) and A3 hyperlinks intact.
foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$4') {
        $h.Delete()
    }
}

# Update the active selection from B7 to B5.
$ws.Range("B5").Select()
